# Update BOM with in-stock parts
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 25 (48MHz Oscillator): swap the out-of-stock Digikey link (535-10086-2-ND)
# for the in-stock replacement part (535-10086-1-ND).
$ws.Range("E25").Value = "https://www.digikey.com/product-detail/en/ASV-48.000MHZ-E-T/535-10086-1-ND/2060881"

# Row 11 (0.1uF / 1206 cap): swap the out-of-stock Kemet Digikey link for the
# in-stock replacement part link.
$ws.Range("E11").Value = "https://www.digikey.com/product-detail/en/C1206F104K1RAC7800/399-5113-1-ND/1465638"

# Move the active selection to F9, matching where the author left off editing.
$ws.Range("F9").Select()
